$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing data rows (2-39): only the columns whose values changed
# (D=Fecha, L=Calidad, M=Volumen, N=Precio minimo, O=Precio maximo,
#  P=Precio promedio ponderado, R=Origen, S=Precio $/Kg)
$updates = @(
    ,@(2, 44165, 'Primera', 300, 27000, 28000, 27500, 'Perú', 1375)
    ,@(3, 44410, 'Primera', 40, 25000, 25000, 25000, 'Perú', 1250)
    ,@(4, 44363, 'Primera', 150, 21000, 22000, 21500, 'Perú', 1075)
    ,@(5, 44442, 'Primera', 30, 22000, 22000, 22000, 'Perú', 1100)
    ,@(6, 44300, 'Primera', 150, 19000, 20000, 19500, 'Perú', 975)
    ,@(7, 44522, 'Primera', 25, 30000, 30000, 30000, 'Perú', 1500)
    ,@(8, 44305, 'Primera', 40, 24000, 24000, 24000, 'Perú', 1200)
    ,@(9, 44620, 'Primera', 60, 22000, 22000, 22000, 'Perú', 1100)
    ,@(10, 44473, 'Primera', 40, 24000, 24000, 24000, 'Perú', 1200)
    ,@(11, 44326, 'Primera', 40, 22000, 22000, 22000, 'Perú', 1100)
    ,@(12, 44372, 'Primera', 60, 20000, 21000, 20667, 'Perú', 1033)
    ,@(13, 44277, 'Primera', 60, 24000, 24000, 24000, 'Perú', 1200)
    ,@(14, 44299, 'Primera', 150, 19000, 20000, 19500, 'Perú', 975)
    ,@(15, 44760, 'Primera', 300, 24000, 25000, 24500, 'Perú', 1225)
    ,@(16, 44166, 'Primera', 120, 28000, 28000, 28000, 'Perú', 1400)
    ,@(17, 44613, 'Primera', 60, 30000, 30000, 30000, 'Perú', 1500)
    ,@(18, 44270, 'Primera', 50, 24000, 24000, 24000, 'Perú', 1200)
    ,@(19, 44350, 'Primera', 90, 21000, 22000, 21556, 'Perú', 1078)
    ,@(20, 44445, 'Primera', 35, 20000, 20000, 20000, 'Perú', 1000)
    ,@(21, 44302, 'Primera', 100, 19000, 20000, 19500, 'Perú', 975)
    ,@(22, 44284, 'Primera', 40, 23000, 23000, 23000, 'Perú', 1150)
    ,@(23, 44312, 'Primera', 50, 22000, 22000, 22000, 'Perú', 1100)
    ,@(24, 44396, 'Primera', 45, 22000, 22000, 22000, 'Perú', 1100)
    ,@(25, 44435, 'Primera', 60, 25000, 25000, 25000, 'Perú', 1250)
    ,@(26, 44333, 'Primera', 30, 22000, 22000, 22000, 'Perú', 1100)
    ,@(27, 44529, 'Primera', 34, 28000, 28000, 28000, 'Perú', 1400)
    ,@(28, 44431, 'Primera', 60, 25000, 25000, 25000, 'Perú', 1250)
    ,@(29, 44263, 'Segunda', 150, 15000, 15000, 15000, 'Perú', 750)
    ,@(30, 44365, 'Primera', 150, 20000, 21000, 20500, 'Perú', 1025)
    ,@(31, 44382, 'Primera', 200, 19000, 20000, 19500, 'Perú', 975)
    ,@(32, 44452, 'Primera', 35, 21000, 22000, 21429, 'Perú', 1071)
    ,@(33, 44417, 'Primera', 30, 24000, 24000, 24000, 'Perú', 1200)
    ,@(34, 44354, 'Primera', 150, 21000, 22000, 21500, 'Perú', 1075)
    ,@(35, 44424, 'Primera', 70, 24000, 25000, 24429, 'Perú', 1221)
    ,@(36, 44356, 'Primera', 100, 20000, 21000, 20500, 'Perú', 1025)
    ,@(37, 44298, 'Primera', 240, 19000, 20000, 19500, 'Perú', 975)
    ,@(38, 44357, 'Primera', 200, 20000, 21000, 20500, 'Perú', 1025)
    ,@(39, 44438, 'Primera', 25, 21000, 21000, 21000, 'Perú', 1050)
)

foreach ($u in $updates) {
    $r = $u[0]
    $ws.Cells.Item($r, 4).Value  = $u[1]   # D Fecha
    $ws.Cells.Item($r, 12).Value = $u[2]   # L Calidad
    $ws.Cells.Item($r, 13).Value = $u[3]   # M Volumen
    $ws.Cells.Item($r, 14).Value = $u[4]   # N Precio minimo
    $ws.Cells.Item($r, 15).Value = $u[5]   # O Precio maximo
    $ws.Cells.Item($r, 16).Value = $u[6]   # P Precio promedio ponderado
    $ws.Cells.Item($r, 18).Value = $u[7]   # R Origen
    $ws.Cells.Item($r, 19).Value = $u[8]   # S Precio $/Kg
}

# Append the new weekly record as row 40
$newRow = @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44355, 13, 'Fruta', 100108, 'Tropicales y subtropicales', 100108007, 'Coco', 'Sin especificar', 'Primera', 200, 20000, 21000, 20500, '$/malla 20 unidades', 'Ecuador', 1025, 20)
for ($c = 0; $c -lt $newRow.Count; $c++) {
    $ws.Cells.Item(40, $c + 1).Value = $newRow[$c]
}

$ws.Range("D40").NumberFormat = $ws.Range("D39").NumberFormat